$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: ERIKA's first logged period corrected to 1902
$ws.Range("E17").Value = "1902"

# Row 18: was a duplicate ERIKA/1911 row -> now the new worker ESTELA PATRICIA POILAO PATRICIA
$ws.Range("C18").Value = "20191183"
$ws.Range("D18").Value = "ESTELA PATRICIA POILAO PATRICIA"
$ws.Range("E18").Value = "1902"
$ws.Range("F18").Value = 32021

# Row 19: was a duplicate ERIKA/1910 row -> now the new worker MAIRA ALEJANDRA JULIO TAPIA
$ws.Range("C19").Value = "20191208"
$ws.Range("D19").Value = "MAIRA ALEJANDRA JULIO TAPIA"
$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 32021

# Row 20: was a duplicate ERIKA/1909 row -> now the new worker ALEJANDRA DEL CARMEN PAYARES CUADRO
$ws.Range("C20").Value = "20191265"
$ws.Range("D20").Value = "ALEJANDRA DEL CARMEN PAYARES CUADRO"
$ws.Range("E20").Value = "1902"

# Row 21: was a duplicate ERIKA/1908 row -> now the new worker OLGA ARAGON MACHUCA
$ws.Range("C21").Value = "20191273"
$ws.Range("D21").Value = "OLGA ARAGON MACHUCA"
$ws.Range("E21").Value = "1902"
$ws.Range("F21").Value = 32021

# Rows 22-27: ERIKA's remaining periods renumbered forward (1907..1902 -> 1903..1908)
$ws.Range("E22").Value = "1903"
$ws.Range("E23").Value = "1904"
$ws.Range("E24").Value = "1905"
$ws.Range("E25").Value = "1906"
$ws.Range("E26").Value = "1907"
$ws.Range("E27").Value = "1908"

# Rows 28-31: were the four new workers (single rows each) -> now ERIKA's remaining periods 1909-1912
$ws.Range("C28").Value = "1047419389"
$ws.Range("D28").Value = "ERIKA PAOLA MARTINEZ SILGADO"
$ws.Range("E28").Value = "1909"
$ws.Range("F28").Value = 33125

$ws.Range("C29").Value = "1047419389"
$ws.Range("D29").Value = "ERIKA PAOLA MARTINEZ SILGADO"
$ws.Range("E29").Value = "1910"
$ws.Range("F29").Value = 33125

$ws.Range("C30").Value = "1047419389"
$ws.Range("D30").Value = "ERIKA PAOLA MARTINEZ SILGADO"
$ws.Range("E30").Value = "1911"

$ws.Range("C31").Value = "1047419389"
$ws.Range("D31").Value = "ERIKA PAOLA MARTINEZ SILGADO"
$ws.Range("E31").Value = "1912"
$ws.Range("F31").Value = 33125
